$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 75 with revised quarterly figures ---
$ws.Range("J75").Value = 5444
$ws.Range("K75").Value = 1605
$ws.Range("L75").Value = 4785
$ws.Range("M75").Value = 6766
$ws.Range("N75").Value = 43123
$ws.Range("O75").Value = 58028
$ws.Range("S75").Value = 19526
$ws.Range("V75").Value = 847
$ws.Range("W75").Value = 29
$ws.Range("X75").Value = 264
$ws.Range("Z75").Value = 6925
$ws.Range("AA75").Value = 1920

# --- Append a new row 76 for period 01-04-2021 ---
# Column A holds the period label as text; force text formatting so the
# "01-04-2021" string isn't auto-converted into a date serial number, then
# restore the default (Normal) cell style so no explicit style is left on
# the cell, matching the rest of the data rows.
$ws.Range("A76").NumberFormat = "@"
$ws.Range("A76").Value = "01-04-2021"
$ws.Range("A76").Style = "Normal"

$ws.Range("B76").Value = 22748
$ws.Range("C76").Value = 22683
$ws.Range("D76").Value = 1
$ws.Range("E76").Value = 64
$ws.Range("F76").Value = 152402
$ws.Range("G76").Value = 46
$ws.Range("H76").Value = 6019
$ws.Range("I76").Value = 15278
$ws.Range("J76").Value = 6272
$ws.Range("K76").Value = 1519
$ws.Range("L76").Value = 4245
$ws.Range("M76").Value = 7603
$ws.Range("N76").Value = 42591
$ws.Range("O76").Value = 66912
$ws.Range("P76").Value = 1915
$ws.Range("Q76").Value = 0
$ws.Range("R76").Value = 0
$ws.Range("S76").Value = 23588
$ws.Range("T76").Value = 2928
$ws.Range("U76").Value = 8215
$ws.Range("V76").Value = 866
$ws.Range("W76").Value = 27
$ws.Range("X76").Value = 345
$ws.Range("Y76").Value = 1425
$ws.Range("Z76").Value = 7790
$ws.Range("AA76").Value = 1991
